# Update TPM-derived values in the NATMI LR-pairs output (Alcam-Chl1)
# This reflects re-computation of ligand/receptor average & total expression
# values (and all values derived from them) using updated TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2057716666666667
$ws.Range("H2").Value = 0.6173149999999999
$ws.Range("I2").Value = 0.01089677771948535
$ws.Range("J2").Value = 0.01089677771948535
$ws.Range("M2").Value = 0.022752
$ws.Range("N2").Value = 0.068256
$ws.Range("O2").Value = 0.003897602285935517
$ws.Range("P2").Value = 0.003897602285935517
$ws.Range("Q2").Value = 0.004681716959999999
$ws.Range("R2").Value = 0.04213545263999999
$ws.Range("S2").Value = 0.0000424713057487973
$ws.Range("T2").Value = 0.0000424713057487973
$ws.Range("G3").Value = 0.2057716666666667
$ws.Range("H3").Value = 0.6173149999999999
$ws.Range("I3").Value = 0.01089677771948535
$ws.Range("J3").Value = 0.01089677771948535
$ws.Range("O3").Value = 0.7259825591205727
$ws.Range("P3").Value = 0.7259825591205727
$ws.Range("Q3").Value = 0.8720348076466669
$ws.Range("R3").Value = 7.84831326882
$ws.Range("S3").Value = 0.007910870574960011
$ws.Range("T3").Value = 0.007910870574960011
$ws.Range("G4").Value = 0.2057716666666667
$ws.Range("H4").Value = 0.6173149999999999
$ws.Range("I4").Value = 0.01089677771948535
$ws.Range("J4").Value = 0.01089677771948535
$ws.Range("M4").Value = 1.576807
$ws.Range("N4").Value = 4.730421
$ws.Range("O4").Value = 0.2701198385934918
$ws.Range("P4").Value = 0.2701198385934918
$ws.Range("Q4").Value = 0.3244622044016666
$ws.Range("R4").Value = 2.920159839615
$ws.Range("S4").Value = 0.002943435838776539
$ws.Range("T4").Value = 0.002943435838776539
$ws.Range("I5").Value = 0.01769706320706529
$ws.Range("J5").Value = 0.01769706320706529
$ws.Range("M5").Value = 0.022752
$ws.Range("N5").Value = 0.068256
$ws.Range("O5").Value = 0.003897602285935517
$ws.Range("P5").Value = 0.003897602285935517
$ws.Range("Q5").Value = 0.007603407455999999
$ws.Range("R5").Value = 0.06843066710399999
$ws.Range("S5").Value = 0.00006897611401020301
$ws.Range("T5").Value = 0.00006897611401020302
$ws.Range("I6").Value = 0.01769706320706529
$ws.Range("J6").Value = 0.01769706320706529
$ws.Range("O6").Value = 0.7259825591205727
$ws.Range("P6").Value = 0.7259825591205727
$ws.Range("S6").Value = 0.01284775923598379
$ws.Range("T6").Value = 0.01284775923598379
$ws.Range("I7").Value = 0.01769706320706529
$ws.Range("J7").Value = 0.01769706320706529
$ws.Range("M7").Value = 1.576807
$ws.Range("N7").Value = 4.730421
$ws.Range("O7").Value = 0.2701198385934918
$ws.Range("P7").Value = 0.2701198385934918
$ws.Range("Q7").Value = 0.5269473497043332
$ws.Range("R7").Value = 4.742526147338999
$ws.Range("S7").Value = 0.004780327857071298
$ws.Range("T7").Value = 0.004780327857071299
$ws.Range("G8").Value = 0.4895776666666666
$ws.Range("H8").Value = 1.468733
$ws.Range("I8").Value = 0.02592591631545138
$ws.Range("J8").Value = 0.02592591631545138
$ws.Range("M8").Value = 0.022752
$ws.Range("N8").Value = 0.068256
$ws.Range("O8").Value = 0.003897602285935517
$ws.Range("P8").Value = 0.003897602285935517
$ws.Range("Q8").Value = 0.011138871072
$ws.Range("R8").Value = 0.100249839648
$ws.Range("S8").Value = 0.0001010489106960762
$ws.Range("T8").Value = 0.0001010489106960762
$ws.Range("G9").Value = 0.4895776666666666
$ws.Range("H9").Value = 1.468733
$ws.Range("I9").Value = 0.02592591631545138
$ws.Range("J9").Value = 0.02592591631545138
$ws.Range("O9").Value = 0.7259825591205727
$ws.Range("P9").Value = 0.7259825591205727
$ws.Range("Q9").Value = 2.074769443702667
$ws.Range("R9").Value = 18.672924993324
$ws.Range("S9").Value = 0.0188217630742372
$ws.Range("T9").Value = 0.0188217630742372
$ws.Range("G10").Value = 0.4895776666666666
$ws.Range("H10").Value = 1.468733
$ws.Range("I10").Value = 0.02592591631545138
$ws.Range("J10").Value = 0.02592591631545138
$ws.Range("M10").Value = 1.576807
$ws.Range("N10").Value = 4.730421
$ws.Range("O10").Value = 0.2701198385934918
$ws.Range("P10").Value = 0.2701198385934918
$ws.Range("Q10").Value = 0.7719694918436666
$ws.Range("R10").Value = 6.947725426592999
$ws.Range("S10").Value = 0.007003104330518103
$ws.Range("T10").Value = 0.007003104330518103
$ws.Range("G11").Value = 17.85418133333333
$ws.Range("H11").Value = 53.562544
$ws.Range("I11").Value = 0.9454802427579979
$ws.Range("J11").Value = 0.945480242757998
$ws.Range("M11").Value = 0.022752
$ws.Range("N11").Value = 0.068256
$ws.Range("O11").Value = 0.003897602285935517
$ws.Range("P11").Value = 0.003897602285935517
$ws.Range("Q11").Value = 0.406218333696
$ws.Range("R11").Value = 3.655965003264
$ws.Range("S11").Value = 0.003685105955480441
$ws.Range("T11").Value = 0.003685105955480441
$ws.Range("G12").Value = 17.85418133333333
$ws.Range("H12").Value = 53.562544
$ws.Range("I12").Value = 0.9454802427579979
$ws.Range("J12").Value = 0.945480242757998
$ws.Range("O12").Value = 0.7259825591205727
$ws.Range("P12").Value = 0.7259825591205727
$ws.Range("Q12").Value = 75.66380657218134
$ws.Range("R12").Value = 680.9742591496321
$ws.Range("S12").Value = 0.6864021662353917
$ws.Range("T12").Value = 0.6864021662353917
$ws.Range("G13").Value = 17.85418133333333
$ws.Range("H13").Value = 53.562544
$ws.Range("I13").Value = 0.9454802427579979
$ws.Range("J13").Value = 0.945480242757998
$ws.Range("M13").Value = 1.576807
$ws.Range("N13").Value = 4.730421
$ws.Range("O13").Value = 0.2701198385934918
$ws.Range("P13").Value = 0.2701198385934918
$ws.Range("Q13").Value = 28.15259810566933
$ws.Range("R13").Value = 253.373382951024
$ws.Range("S13").Value = 0.2553929705671258
$ws.Range("T13").Value = 0.2553929705671259

Write-Output "Updated TPM-derived values for Alcam-Chl1 sheet"
